$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B-G (ME, MAE, MSE, RMSE, SE, N)
$data = @{
    2  = @(0.5692650327263563, 0.5692650327263563, 0.3908711553892542, 0.6251968932978268, 0.2682304822149452, 14)
    3  = @(0.3997456872501602, 0.4033467742491634, 0.2178465991588477, 0.4667403980360472, 0.2507737693766841, 13)
    4  = @(0.3400043916411057, 0.3606492384221238, 0.1784350047575419, 0.4224156776891003, 0.261809129132386, 12)
    5  = @(0.4288908803047028, 0.4288908803047028, 0.2376241711832545, 0.4874670975391617, 0.242990663137872, 11)
    6  = @(0.3767221034172891, 0.3801506752190701, 0.1848819966348901, 0.429979065344919, 0.2184858333361773, 10)
    7  = @(0.3516807317407905, 0.358904738120446, 0.1712159718543759, 0.4137825175794355, 0.2312546521125961, 9)
    8  = @(0.3826977698461033, 0.3826977698461033, 0.193400430583437, 0.4397731580979414, 0.2316224700134891, 8)
    9  = @(0.3590813353308283, 0.3590813353308283, 0.1682323819789578, 0.4101614096657044, 0.2141069966830871, 7)
    10 = @(0.3470367157308435, 0.3470367157308435, 0.1621328538698305, 0.4026572411739673, 0.2236918553848042, 6)
    11 = @(0.308333697360539, 0.308333697360539, 0.1232078055746762, 0.3510096944169437, 0.1875437837101504, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
